$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21, pushing the existing rows 21-29 down to 22-30.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Vega Monumental Concepción"
$ws.Range("C21").Value = "Bíobío"
$ws.Range("D21").Value = "2022-04-20"
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 100112022
$ws.Range("G21").Value = "Arveja Verde"
$ws.Range("H21").Value = "Perfection"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 110
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 26000
$ws.Range("M21").Value = 25545
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Carahue"
$ws.Range("P21").Value = 1022
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
